$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds quarterly dates (1st day of the quarter-start month: Jan/Apr/Jul/Oct).
# Bug fix: shift each date to the 15th of the following month.
for ($r = 2; $r -le 150; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $serial = $cellA.Value2
    $d = [DateTime]::FromOADate($serial)
    $d2 = $d.AddMonths(1)
    $d3 = $d2.AddDays(15 - $d2.Day)
    $cellA.Value = $d3
}
